$wb = $excel.ActiveWorkbook

# Add the new worksheet and move it to the end of the workbook (after the
# existing "Testcase1" sheet) so it becomes sheet 2
$newSheet = $wb.Worksheets.Add()
$newSheet.Name = "VerifyActiTimeLogin"
$newSheet.Move($null, $wb.Worksheets.Item($wb.Worksheets.Count))

# Sheet objects are tracked by position, and Move() re-indexes the
# collection, so re-fetch the handle by name to keep referring to the
# newly added sheet
$newSheet = $wb.Worksheets.Item("VerifyActiTimeLogin")

# Populate the login test-data table
$newSheet.Range("A1").Value = "UserName"
$newSheet.Range("B1").Value = "Password"
$newSheet.Range("A2").Value = "admin"
$newSheet.Range("B2").Value = "manager"
$newSheet.Range("A3").Value = "pawan"
$newSheet.Range("B3").Value = "manager"

# Select the final data cell and make this the active sheet/tab
$newSheet.Range("B3").Select()
$newSheet.Activate()
